# Update info clase 04 File System
# Adds a new "Clase 04" section (header row + 18 video-index rows) to the
# index sheet, reusing the existing header / data row formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header row (45): same look as the existing "Clase NN" rows ---
$ws.Range("A8:B8").Copy()
$ws.Range("A45:B45").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A45").Value = "Clase 04"

# --- Data rows (46-64): same look as the existing time/description rows ---
$ws.Range("A2:B2").Copy()
$ws.Range("A46:B64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A46").Value = 0.003125
$ws.Range("B46").Value = "Que es Node JS"

$ws.Range("A47").Value = 0.004861111111111111
$ws.Range("B47").Value = "Actividad en clase: Proyecto de node (generar 10000 numeros aleatorios de 1 a 20)"

$ws.Range("A48").Value = 0.008680555555555556
$ws.Range("B48").Value = "Uso de metodo Math"

$ws.Range("A49").Value = 0.019791666666666666
$ws.Range("B49").Value = "Modulos nativos de Nodejs"

$ws.Range("A50").Value = 0.021180555555555557
$ws.Range("B50").Value = "Observacion en cambio del orden de las clase, primero fs despues crypto"

$ws.Range("A51").Value = 0.028819444444444446
$ws.Range("B51").Value = "Ejemplo de setTimeout"

$ws.Range("A52").Value = 0.03090277777777778
$ws.Range("B52").Value = "Ejemplo de setTimeinterval"

$ws.Range("A53").Value = 0.03194444444444444
$ws.Range("B53").Value = "Manejo de archivos - La persistencia en memoria"

$ws.Range("A54").Value = 0.034027777777777775
$ws.Range("B54").Value = "file system en NodeJs"

$ws.Range("A55").Value = 0.0375
$ws.Range("B55").Value = "metodo writeFileSync()"

$ws.Range("A56").Value = 0.03819444444444445
$ws.Range("B56").Value = "metodo appendFileSync()"

$ws.Range("A57").Value = 0.03888888888888889
$ws.Range("B57").Value = "metodo existsSync()"

# NOTE: shared-string pool order follows write order, not row order - the
# source workbook registered "metodo readFile()" (row 59) before "File
# System de manera Asincronica..." (row 58), so the B59 value is written
# first to reproduce that exact sharedStrings.xml ordering.
$ws.Range("B59").Value = "metodo readFile()"
$ws.Range("A59").Value = 0.044444444444444446

$ws.Range("A58").Value = 0.043402777777777776
$ws.Range("B58").Value = "File System de manera Asincronica (promises) - no tienen --Sync-- al final"

$ws.Range("A60").Value = 0.04548611111111111
$ws.Range("B60").Value = "metodo appendFile()"

$ws.Range("A61").Value = 0.05173611111111111
$ws.Range("B61").Value = "archivos .json"

$ws.Range("A62").Value = 0.052083333333333336
$ws.Range("B62").Value = 'metodo JSON.stringify() - para pasar de tipo objeto JS a texto plano de tipo JSON '

$ws.Range("A63").Value = 0.053125
$ws.Range("B63").Value = 'metodo JSON.parse() - para pasar de texto plano de tipo JSON  a tipo objeto JS'

$ws.Range("A64").Value = 0.06145833333333333
$ws.Range("B64").Value = 'Hands on Labs - "Manager de usuarios"'

# --- Bring the view in line with where the author ended up editing ---
$ws.Range("A65").Select() | Out-Null
